$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.524.58"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.893.14"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'237.78"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.4901"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").Value = "'0.2935"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").Value = "1.891.10"
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").Value = "'17.12"
$ws.Range("E11").Value = "  +2.78%  "
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").Value = "'5.162"
$ws.Range("E13").Value = "  +3.14%  "
$ws.Range("D14").Value = "'88.16"
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").Value = "'0.6677"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Value = "30.471.76"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "'13.41"
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("D18").Value = "'0.000007833"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "2.163.93"
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("D21").Value = "'5.319"
$ws.Range("E21").Value = "  +12.34%  "
$ws.Range("D22").Value = "'1.004"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "'191.72"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").Value = "'6.181"
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("D25").Value = "'9.487"
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("D26").Value = "'163.16"
$ws.Range("E26").Value = "  +3.21%  "
$ws.Range("D27").Value = "'18.37"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").Value = "'1.937"
$ws.Range("E28").Value = "  +5.69%  "
$ws.Range("D29").Value = "'1.472"
$ws.Range("E29").Value = "  +4.99%  "
$ws.Range("D30").Value = "'4.381"
$ws.Range("E30").Value = "  +2.73%  "
$ws.Range("D31").Value = "'0.09158"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("D32").Value = "'4.099"
$ws.Range("E32").Value = "  +3.98%  "
$ws.Range("D33").Value = "'0.05219"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").Value = "'0.7404"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").Value = "'1.099"
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("D36").Value = "'2.719"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").Value = "'0.01829"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").Value = "'2.684"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").Value = "'0.9187"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").Value = "'2.050"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").Value = "'0.4406"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "'5.957"
$ws.Range("E42").Value = "  +3.95%  "
$ws.Range("D43").Value = "'106.23"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").Value = "'0.9937"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("E45").Value = "  +3.07%  "
$ws.Range("D46").Value = "'68.79"
$ws.Range("E46").Value = "  +20.20%  "
$ws.Range("D47").Value = "'7.643"
$ws.Range("E47").Value = "  +4.79%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.025"
$ws.Range("E48").Value = "  +3.99%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'34.97"
$ws.Range("E49").Value = "  +5.25%  "
$ws.Range("D50").Value = "'0.05831"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("E51").Value = "  -3.74%  "
